$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Props")

# Update the "onKeyDown" row (row 18) destinations: add "2 - Root" and
# "3 - ScrollingContainer" as its destinations.
$ws.Range("C18").Value = "2 - Root"
$ws.Range("D18").Value = "3 - ScrollingContainer"

# Remove the "liveColumnResize" row entirely (row 19); this shifts every
# row below it up by one.
$ws.Rows.Item(19).Delete()
